# Applies the "Office schema" commit:
#   - Sheet "Normalization exercersie 1_31Ja": adds a new "Final table
#     list in 3NF:" summary block (rows 48-59) that re-lists the three
#     final 3NF tables (Students / Course / Faculty) and their columns.
#   - Sheet "Normalization exercersie 2_31Ja":
#       * fixes two cells that referenced a stray typo'd string "Brime"
#         so they read "Btime" instead (R27, S38);
#       * renames the "emp"/"dept" table-name cells to "employee" /
#         "department" (G37, K37);
#       * adds the matching new "Final table list in 3NF:" summary block
#         (rows 47-58) re-listing the four final 3NF tables (employee /
#         department / project / btime) and their columns.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

function Set-FormattedValue {
    param($ws, [string]$target, [string]$sourceForFormat, $value)

    $ws.Range($sourceForFormat).Copy()
    $ws.Range($target).PasteSpecial($xlPasteFormats)
    if ($null -ne $value) {
        $ws.Range($target).Value = $value
    }
}

# ---------------------------------------------------------------------
# Sheet 1: "Normalization exercersie 1_31Ja"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Normalization exercersie 1_31Ja")

# Section title (same style as the other section headers, e.g. B37).
Set-FormattedValue $ws1 "B48" "B37" "Final table list in 3NF:"

# Header filler row (mirrors F37/F38:H38's odd reuse of the "Grade" text).
Set-FormattedValue $ws1 "B49" "F37" "Grade"
Set-FormattedValue $ws1 "C49" "K38" $null
Set-FormattedValue $ws1 "D49" "K38" $null

Set-FormattedValue $ws1 "B50" "F38" "SID*"
Set-FormattedValue $ws1 "C50" "G38" "CID*"
Set-FormattedValue $ws1 "D50" "H38" "Grade"

# Students table
Set-FormattedValue $ws1 "B52" "J38" "Students table"
Set-FormattedValue $ws1 "C52" "K38" $null
Set-FormattedValue $ws1 "B53" "J39" "SID"
Set-FormattedValue $ws1 "C53" "K39" "S_name"

# Course table
Set-FormattedValue $ws1 "B55" "M38" "Course table"
Set-FormattedValue $ws1 "C55" "K38" $null
Set-FormattedValue $ws1 "B56" "M39" "CID"
Set-FormattedValue $ws1 "C56" "N39" "FID*"
Set-FormattedValue $ws1 "D56" "O39" "C_name"

# Faculty table
Set-FormattedValue $ws1 "B58" "Q38" "Faculty table"
Set-FormattedValue $ws1 "C58" "K38" $null
Set-FormattedValue $ws1 "B59" "Q39" "FID"
Set-FormattedValue $ws1 "C59" "R39" "F_name"
Set-FormattedValue $ws1 "D59" "S39" "F_phone"

# ---------------------------------------------------------------------
# Sheet 2: "Normalization exercersie 2_31Ja"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Normalization exercersie 2_31Ja")

# Typo fix: the header cells pointed at a duplicated "Brime" string;
# point them at the correctly spelled "Btime" instead (same text the
# rest of the workbook already uses for this column).
$ws2.Range("R27").Value = "Btime"
$ws2.Range("S38").Value = "Btime"

# Rename the abbreviated table names to their full forms.
$ws2.Range("G37").Value = "employee"
$ws2.Range("K37").Value = "department"

# Section title (same style as the other section headers, e.g. B37).
Set-FormattedValue $ws2 "B47" "B37" "Final table list in 3NF:"

# employee table
Set-FormattedValue $ws2 "B48" "G37" "employee"
Set-FormattedValue $ws2 "C48" "K27" $null
Set-FormattedValue $ws2 "D48" "K27" $null
Set-FormattedValue $ws2 "B49" "G38" "EID"
Set-FormattedValue $ws2 "C49" "H38" "Ename"
Set-FormattedValue $ws2 "D49" "I38" "DID*"

# department table
Set-FormattedValue $ws2 "B51" "K37" "department"
Set-FormattedValue $ws2 "C51" "K27" $null
Set-FormattedValue $ws2 "B52" "K38" "DID"
Set-FormattedValue $ws2 "C52" "L38" "Dname"

# project table
Set-FormattedValue $ws2 "B54" "N37" "project"
Set-FormattedValue $ws2 "C54" "K27" $null
Set-FormattedValue $ws2 "B55" "N38" "PID"
Set-FormattedValue $ws2 "C55" "O38" "Pname"

# btime (junction) table
Set-FormattedValue $ws2 "B57" "Q37" "btime"
Set-FormattedValue $ws2 "C57" "K27" $null
Set-FormattedValue $ws2 "D57" "K27" $null
Set-FormattedValue $ws2 "B58" "Q38" "EID*"
Set-FormattedValue $ws2 "C58" "R38" "PID*"
Set-FormattedValue $ws2 "D58" "S38" "Btime"
